# Workbook currently has one sheet ("Sheet") with "Hello Excel" in A1.
# This script:
#   1. Adds two more cells of text to Sheet1 (A2, K11) -- growing the
#      used range to A1:K11.
#   2. Inserts a brand-new worksheet "Sheet2" right after Sheet1, with a
#      single greeting in A1 and the selection parked on C3.
#   3. Leaves Sheet1 as the active/selected tab, cursor still on A1.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# --- Sheet1: extra cells ---------------------------------------------
$ws1.Range("A2").Value  = "lecxE olleH"
$ws1.Range("K11").Value = "Greetings"

# --- Sheet2: new worksheet, placed after Sheet1 -----------------------
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Sheet2"
$ws2.Range("A1").Value = "Hi Olivia"
$ws2.Range("C3").Select() | Out-Null

# Re-activate Sheet1 so it remains the selected/visible tab.
$ws1.Select() | Out-Null
